$wb = $excel.ActiveWorkbook

# --- Rename header cells (B1) on the existing sheets ---
$wsWeekly = $wb.Worksheets.Item("Weekly Quantity")
$wsWeekly.Range("B1").Value = "Weekly_PO_Qty"

$wsMonthly = $wb.Worksheets.Item("Monthly Trend")
$wsMonthly.Range("B1").Value = "Monthly_PO_Qty"

# --- Add the new "PO Forecast" sheet after the last existing sheet ---
$sheetCount = $wb.Worksheets.Count
$lastSheet = $wb.Worksheets.Item($sheetCount)
$wsForecast = $wb.Worksheets.Add($null, $lastSheet)
$wsForecast.Name = "PO Forecast"

# --- Header row ---
$wsForecast.Range("A1").Value = "ds"
$wsForecast.Range("B1").Value = "PO_Forecast"
$wsForecast.Range("C1").Value = "yhat_lower"
$wsForecast.Range("D1").Value = "yhat_upper"

# Reuse the bold/bordered header style already present in the workbook
# (copy format only, so no new style entries are introduced).
$wsWeekly.Range("A1:B1").Copy()
$wsForecast.Range("A1:B1").PasteSpecial(-4122)
$wsWeekly.Range("A1:B1").Copy()
$wsForecast.Range("C1:D1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$wsForecast.Cells.Item(2,1).Value = 45018.99999999999
$wsForecast.Cells.Item(2,2).Value = 2
$wsForecast.Cells.Item(2,3).Value = -6.868350049772292
$wsForecast.Cells.Item(2,4).Value = 10.81995647660504
$wsForecast.Cells.Item(3,1).Value = 45039.99999999999
$wsForecast.Cells.Item(3,2).Value = 2
$wsForecast.Cells.Item(3,3).Value = -7.473894049938778
$wsForecast.Cells.Item(3,4).Value = 11.23496616969146
$wsForecast.Cells.Item(4,1).Value = 45046.99999999999
$wsForecast.Cells.Item(4,2).Value = 2
$wsForecast.Cells.Item(4,3).Value = -6.922534580808817
$wsForecast.Cells.Item(4,4).Value = 10.61782062776201
$wsForecast.Cells.Item(5,1).Value = 45067.99999999999
$wsForecast.Cells.Item(5,2).Value = 2
$wsForecast.Cells.Item(5,3).Value = -6.253589088158839
$wsForecast.Cells.Item(5,4).Value = 11.21571330609387
$wsForecast.Cells.Item(6,1).Value = 45081.99999999999
$wsForecast.Cells.Item(6,2).Value = 3
$wsForecast.Cells.Item(6,3).Value = -5.967788177755618
$wsForecast.Cells.Item(6,4).Value = 11.60684181486866
$wsForecast.Cells.Item(7,1).Value = 45088.99999999999
$wsForecast.Cells.Item(7,2).Value = 3
$wsForecast.Cells.Item(7,3).Value = -5.527661507522686
$wsForecast.Cells.Item(7,4).Value = 11.62174713437897
$wsForecast.Cells.Item(8,1).Value = 45116.99999999999
$wsForecast.Cells.Item(8,2).Value = 3
$wsForecast.Cells.Item(8,3).Value = -5.150683649404727
$wsForecast.Cells.Item(8,4).Value = 11.53606957152183
$wsForecast.Cells.Item(9,1).Value = 45123.99999999999
$wsForecast.Cells.Item(9,2).Value = 3
$wsForecast.Cells.Item(9,3).Value = -5.893404444517508
$wsForecast.Cells.Item(9,4).Value = 12.87533569331439
$wsForecast.Cells.Item(10,1).Value = 45130.99999999999
$wsForecast.Cells.Item(10,2).Value = 3
$wsForecast.Cells.Item(10,3).Value = -5.91887224943324
$wsForecast.Cells.Item(10,4).Value = 11.75168830602704
$wsForecast.Cells.Item(11,1).Value = 45137.99999999999
$wsForecast.Cells.Item(11,2).Value = 4
$wsForecast.Cells.Item(11,3).Value = -5.038762047104713
$wsForecast.Cells.Item(11,4).Value = 12.44517620457103
$wsForecast.Cells.Item(12,1).Value = 45144.99999999999
$wsForecast.Cells.Item(12,2).Value = 4
$wsForecast.Cells.Item(12,3).Value = -5.292372574998972
$wsForecast.Cells.Item(12,4).Value = 12.33585659282154
$wsForecast.Cells.Item(13,1).Value = 45151.99999999999
$wsForecast.Cells.Item(13,2).Value = 4
$wsForecast.Cells.Item(13,3).Value = -4.789949332292152
$wsForecast.Cells.Item(13,4).Value = 12.88403151592122
$wsForecast.Cells.Item(14,1).Value = 45158.99999999999
$wsForecast.Cells.Item(14,2).Value = 4
$wsForecast.Cells.Item(14,3).Value = -5.203524533410732
$wsForecast.Cells.Item(14,4).Value = 12.37920719294642
$wsForecast.Cells.Item(15,1).Value = 45165.99999999999
$wsForecast.Cells.Item(15,2).Value = 4
$wsForecast.Cells.Item(15,3).Value = -5.094013217780366
$wsForecast.Cells.Item(15,4).Value = 13.07227844509258
$wsForecast.Cells.Item(16,1).Value = 45172.99999999999
$wsForecast.Cells.Item(16,2).Value = 4
$wsForecast.Cells.Item(16,3).Value = -4.516939340833722
$wsForecast.Cells.Item(16,4).Value = 13.1537058567212
$wsForecast.Cells.Item(17,1).Value = 45186.99999999999
$wsForecast.Cells.Item(17,2).Value = 4
$wsForecast.Cells.Item(17,3).Value = -4.554481011131836
$wsForecast.Cells.Item(17,4).Value = 13.31016863252553
$wsForecast.Cells.Item(18,1).Value = 45445.99999999999
$wsForecast.Cells.Item(18,2).Value = 9
$wsForecast.Cells.Item(18,3).Value = -0.2969641637363782
$wsForecast.Cells.Item(18,4).Value = 17.40059522990411
$wsForecast.Cells.Item(19,1).Value = 45459.99999999999
$wsForecast.Cells.Item(19,2).Value = 9
$wsForecast.Cells.Item(19,3).Value = -0.1732228681901788
$wsForecast.Cells.Item(19,4).Value = 17.50644227940835
$wsForecast.Cells.Item(20,1).Value = 45480.99999999999
$wsForecast.Cells.Item(20,2).Value = 9
$wsForecast.Cells.Item(20,3).Value = 0.2083616301802291
$wsForecast.Cells.Item(20,4).Value = 18.40093102456071
$wsForecast.Cells.Item(21,1).Value = 45487.99999999999
$wsForecast.Cells.Item(21,2).Value = 9
$wsForecast.Cells.Item(21,3).Value = 0.2096648565571771
$wsForecast.Cells.Item(21,4).Value = 18.11723861962371
$wsForecast.Cells.Item(22,1).Value = 45494.99999999999
$wsForecast.Cells.Item(22,2).Value = 9
$wsForecast.Cells.Item(22,3).Value = 0.1059668511391858
$wsForecast.Cells.Item(22,4).Value = 18.67351473605877
$wsForecast.Cells.Item(23,1).Value = 45515.99999999999
$wsForecast.Cells.Item(23,2).Value = 10
$wsForecast.Cells.Item(23,3).Value = 1.202091454466805
$wsForecast.Cells.Item(23,4).Value = 18.21914238015036
$wsForecast.Cells.Item(24,1).Value = 45550.99999999999
$wsForecast.Cells.Item(24,2).Value = 10
$wsForecast.Cells.Item(24,3).Value = 0.9381433697250025
$wsForecast.Cells.Item(24,4).Value = 19.56552251858779
$wsForecast.Cells.Item(25,1).Value = 45620.99999999999
$wsForecast.Cells.Item(25,2).Value = 11
$wsForecast.Cells.Item(25,3).Value = 2.832065530025953
$wsForecast.Cells.Item(25,4).Value = 20.19399182552035
$wsForecast.Cells.Item(26,1).Value = 45634.99999999999
$wsForecast.Cells.Item(26,2).Value = 12
$wsForecast.Cells.Item(26,3).Value = 2.500195316556124
$wsForecast.Cells.Item(26,4).Value = 20.69091838556769
$wsForecast.Cells.Item(27,1).Value = 45641.99999999999
$wsForecast.Cells.Item(27,2).Value = 12
$wsForecast.Cells.Item(27,3).Value = 2.700989854420117
$wsForecast.Cells.Item(27,4).Value = 21.0177566696163
$wsForecast.Cells.Item(28,1).Value = 45648.99999999999
$wsForecast.Cells.Item(28,2).Value = 12
$wsForecast.Cells.Item(28,3).Value = 3.449372424108407
$wsForecast.Cells.Item(28,4).Value = 19.81038376188928
$wsForecast.Cells.Item(29,1).Value = 45655.99999999999
$wsForecast.Cells.Item(29,2).Value = 12
$wsForecast.Cells.Item(29,3).Value = 2.923236893740314
$wsForecast.Cells.Item(29,4).Value = 20.35834621862873
$wsForecast.Cells.Item(30,1).Value = 45662.99999999999
$wsForecast.Cells.Item(30,2).Value = 12
$wsForecast.Cells.Item(30,3).Value = 3.422470700949631
$wsForecast.Cells.Item(30,4).Value = 20.77366654853027
$wsForecast.Cells.Item(31,1).Value = 45669.99999999999
$wsForecast.Cells.Item(31,2).Value = 12
$wsForecast.Cells.Item(31,3).Value = 3.158284442832333
$wsForecast.Cells.Item(31,4).Value = 20.63092036930002
$wsForecast.Cells.Item(32,1).Value = 45676.99999999999
$wsForecast.Cells.Item(32,2).Value = 12
$wsForecast.Cells.Item(32,3).Value = 3.735247154892078
$wsForecast.Cells.Item(32,4).Value = 21.08852343426644
$wsForecast.Cells.Item(33,1).Value = 45683.99999999999
$wsForecast.Cells.Item(33,2).Value = 12
$wsForecast.Cells.Item(33,3).Value = 3.692891404828732
$wsForecast.Cells.Item(33,4).Value = 20.68746201161107
$wsForecast.Cells.Item(34,1).Value = 45690.99999999999
$wsForecast.Cells.Item(34,2).Value = 12
$wsForecast.Cells.Item(34,3).Value = 3.413447674892486
$wsForecast.Cells.Item(34,4).Value = 21.63964096558609
$wsForecast.Cells.Item(35,1).Value = 45697.99999999999
$wsForecast.Cells.Item(35,2).Value = 13
$wsForecast.Cells.Item(35,3).Value = 4.043286155870315
$wsForecast.Cells.Item(35,4).Value = 21.16417918359964

# Apply the same date/time number-format style used for column A on the
# other sheets to the new "ds" column (reuses the existing style, too).
$wsWeekly.Range("A2").Copy()
$wsForecast.Range("A2:A35").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Match the page margins used elsewhere in the workbook (0.75in / 1in / 0.5in).
$wsForecast.PageSetup.LeftMargin = 54
$wsForecast.PageSetup.RightMargin = 54
$wsForecast.PageSetup.TopMargin = 72
$wsForecast.PageSetup.BottomMargin = 72
$wsForecast.PageSetup.HeaderMargin = 36
$wsForecast.PageSetup.FooterMargin = 36

# Leave the workbook selection on the first sheet, as before the edit.
[void]$wsWeekly.Select()
[void]$wsWeekly.Range("A1").Select()
